$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("F10").Value = "https://www.engineersedge.com/properties_of_metals.htm#:~:text=Thermal%20Properties%20of%20Metals%2C%20Conductivity%2C%20Thermal%20Expansion%2C%20Specific%20Heat,-Heat%20Transfer%20Engineering&text=Metals%20in%20general%20have%20high,metals%20are%20shiny%20and%20lustrous"
